$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: InsertXML payload wrapper (OOXML "flat OPC" package fragment)
# ---------------------------------------------------------------------------
function New-PackageXml([string]$bodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + $bodyXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------------
# 1) "Multiples of 2 loop" -> 3 runs: "Multiples of " / "3" / " loop"
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Multiples of 2 loop", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $target = $d.Range($rng.Start, $rng.End + 1)
    $body = '<w:body><w:p w14:paraId="6E054493" w14:textId="0677866B" w:rsidR="00231FDC" w:rsidRDefault="00231FDC" w:rsidP="00A626AD"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="18"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Multiples of </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>3</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> loop</w:t></w:r></w:p><w:sectPr/></w:body>'
    $target.InsertXML((New-PackageXml $body))
}

# ---------------------------------------------------------------------------
# 2) "Looped Number Guessing Game" -> "Palindrome Checker" (simple text swap,
#    single run, formatting unaffected)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Looped Number Guessing Game", $true, $false, $false, $false, $false, $true, 1, $false, "Palindrome Checker", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "Imports and JSON parsing" -> 3 runs: "Imports and JSON" / "/TXT" / " parsing"
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Imports and JSON parsing", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $target = $d.Range($rng.Start, $rng.End + 1)
    $body = '<w:body><w:p w14:paraId="5EA4EDDD" w14:textId="3FD6EFDC" w:rsidR="00231FDC" w:rsidRDefault="00231FDC" w:rsidP="00DC057B"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="24"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Imports and JSON</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>/TXT</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> parsing</w:t></w:r></w:p><w:sectPr/></w:body>'
    $target.InsertXML((New-PackageXml $body))
}

# ---------------------------------------------------------------------------
# 4) Empty "Take Home Challenge" paragraph after "Store in JSON" gains a
#    numPr (numId 27) and the text "Word Counter" + " from TXT file"
#    (two runs). The paragraph is the very last one in its table cell, so it
#    has to be populated in two steps: first insert plain text into it via
#    the paragraph's own Range (collapsed-range inserts at this position are
#    unreliable), then re-locate it by its new text and rewrite the whole
#    paragraph (pPr + runs) via InsertXML now that it is addressable like an
#    ordinary paragraph.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Store in JSON", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $storyEnd = $rng.End
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Start -eq $storyEnd + 1) {
            $p.Range.InsertBefore("Word Counter from TXT file")
            break
        }
    }

    $rng2 = $d.Content
    $found2 = $rng2.Find.Execute("Word Counter from TXT file", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found2) {
        $target2 = $d.Range($rng2.Start, $rng2.End + 1)
        $body2 = '<w:body><w:p w14:paraId="05BF876E" w14:textId="391035D1" w:rsidR="00FA58F3" w:rsidRPr="00CA6BEF" w:rsidRDefault="00FA58F3" w:rsidP="00FA58F3"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="27"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Word Counter</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> from TXT file</w:t></w:r></w:p><w:sectPr/></w:body>'
        $target2.InsertXML((New-PackageXml $body2))
    }
}
